# Update the division problems in the table to the new set of values.
# Cells are addressed by (row, column) so that the two cells that happen
# to share the same original text ("93÷5=") are updated independently
# and correctly, each receiving its own replacement value.

$d = $word.ActiveDocument
$t = $d.Tables(1)

$t.Cell(1, 1).Range.Text  = "20÷5="   # was 93÷5=
$t.Cell(1, 2).Range.Text  = "11÷2="   # was 73÷7=
$t.Cell(1, 3).Range.Text  = "96÷3="   # was 59÷8=
$t.Cell(1, 4).Range.Text  = "15÷3="   # was 78÷8=
$t.Cell(1, 5).Range.Text  = "25÷7="   # was 72÷8=

$t.Cell(5, 1).Range.Text  = "13÷6="   # was 15÷8=
$t.Cell(5, 2).Range.Text  = "29÷4="   # was 91÷7=
$t.Cell(5, 3).Range.Text  = "48÷6="   # was 25÷2=
$t.Cell(5, 4).Range.Text  = "71÷4="   # was 23÷2=
$t.Cell(5, 5).Range.Text  = "26÷7="   # was 75÷2=

$t.Cell(9, 1).Range.Text  = "88÷7="   # was 12÷2=
$t.Cell(9, 2).Range.Text  = "84÷4="   # was 24÷8=
$t.Cell(9, 3).Range.Text  = "63÷2="   # was 59÷3=
$t.Cell(9, 4).Range.Text  = "29÷7="   # was 36÷9=
$t.Cell(9, 5).Range.Text  = "11÷7="   # was 37÷8=

$t.Cell(13, 1).Range.Text = "57÷3="   # was 90÷9=
$t.Cell(13, 2).Range.Text = "63÷8="   # was 96÷4=
$t.Cell(13, 3).Range.Text = "34÷3="   # was 96÷5=
$t.Cell(13, 4).Range.Text = "66÷4="   # was 43÷8=
$t.Cell(13, 5).Range.Text = "87÷9="   # was 55÷6=

$t.Cell(17, 1).Range.Text = "17÷2="   # was 93÷7=
$t.Cell(17, 2).Range.Text = "98÷2="   # was 76÷4=
$t.Cell(17, 3).Range.Text = "91÷5="   # was 93÷5=
$t.Cell(17, 4).Range.Text = "36÷2="   # was 85÷4=
$t.Cell(17, 5).Range.Text = "80÷9="   # was 34÷2=
